# Actualización automática 2025-06-16 13:01:14
#
# Sheet "VENTAS POR GRUPO" gets a new "GRANITO" column inserted right
# before "GRIFERIAS" (i.e. at column F), shifting the existing
# GRIFERIAS..SAL SOLUBLE columns one position to the right, and three
# brand new columns (NO RESURTIBLES, PANELES PVC, PANELES PU) appended
# after the (now shifted) SAL SOLUBLE column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Excel's ColumnWidth COM property is offset from the stored OOXML
# column width by this constant (empirically: stored = ColumnWidth + offset).
$widthOffset = 0.8333333333333334

$lastDataRow = 57    # rows 2..57 hold the per-client numeric data

# ---------------------------------------------------------------
# 1. Insert a brand new column at F ("GRANITO"). This shifts the old
#    F..N columns (GRIFERIAS .. SAL SOLUBLE) one place to the right,
#    becoming G..O, carrying along their values/styles automatically.
#    The inserted column inherits formatting from its left neighbor,
#    which already matches what we need (header style, data style,
#    summary-row style).
# ---------------------------------------------------------------
$ws.Columns.Item(6).Insert()

# Header
$ws.Range("F1").Value = "GRANITO"

# Data rows: new column defaults to 0
for ($r = 2; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# Summary row 58
$ws.Range("F58").Value = "0 de 56"

# Restore the column width for F (GRANITO) -> stored width 13
$ws.Columns.Item(6).ColumnWidth = 13 - $widthOffset

# ---------------------------------------------------------------
# 2. Append three new columns at the end: P (NO RESURTIBLES),
#    Q (PANELES PVC), R (PANELES PU). Copy formatting from the
#    (now shifted) O column ("SAL SOLUBLE") so headers/data/summary
#    cells reuse the same styles as the rest of the table.
# ---------------------------------------------------------------
$ws.Range("O1").EntireColumn.Copy()
$ws.Range("P1:R1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Headers
$ws.Range("P1").Value = "NO RESURTIBLES"
$ws.Range("Q1").Value = "PANELES PVC"
$ws.Range("R1").Value = "PANELES PU"

# Data rows default to 0
for ($r = 2; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0
    $ws.Cells.Item($r, 18).Value = 0
}

# Row 50 (SARITAMA HERRERA MARIA ELIZABETH) has a genuine NO RESURTIBLES value
$ws.Range("P50").Value = 9.58

# Summary row 58
$ws.Range("P58").Value = "1 de 56"
$ws.Range("Q58").Value = "0 de 56"
$ws.Range("R58").Value = "0 de 56"

# Column widths: P=20, Q=17, R=16
$ws.Columns.Item(16).ColumnWidth = 20 - $widthOffset
$ws.Columns.Item(17).ColumnWidth = 17 - $widthOffset
$ws.Columns.Item(18).ColumnWidth = 16 - $widthOffset

# Dimension becomes A1:R58 automatically as Excel recalculates the used range.
